$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove F4 content (shared string "3d") - clear the cell
$ws.Range("F4").ClearContents()

# Update selection to H11 as per the diff
$ws.Range("H11").Select()
